# process.xlsx — "complete loainhanvien & khachhang"
#
# Updates the "Thực tế" (actual status) column on the weekly plan sheet:
#   - Row 13 (Phân quyền truy cập cơ sở dữ liệu)            -> hoàn thành
#   - Row 14 (Tạo form nhân viên, khách hàng, loại nhân viên) -> hoàn thành
#   - Row 22 (Giao hàng...)                                   -> Chưa làm kịp
# and moves the active selection to E23.
#
# NOTE on write order: F22 is written before F13/F14 so that the new shared
# strings are appended to sharedStrings.xml in the same order the source
# workbook uses them ("Chưa làm kịp " before "hoàn thành").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F22").Value = "Chưa làm kịp "
$ws.Range("F13").Value = "hoàn thành"
$ws.Range("F14").Value = "hoàn thành"

$ws.Range("E23").Select()
